$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to write a value as literal TEXT (preserving the cell's existing
# NumberFormat / style index) instead of letting Excel auto-coerce a
# numeric-looking string into a real number.
function Set-TextValue($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $origFmt = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.NumberFormat = $origFmt
}

# ------------------------------------------------------------------
# 1) Insert a new product row before row 13 (pushes DICLOPRO.. down by
#    one row, copying formatting/merges style from the row above).
# ------------------------------------------------------------------
$ws.Rows.Item(13).Insert()

# Re-create the merged cells for the freshly inserted row 13 (Insert()
# does not automatically re-merge the new row).
$ws.Range("A13:B13").Merge()
$ws.Range("C13:G13").Merge()
$ws.Range("H13:K13").Merge()
$ws.Range("L13:M13").Merge()
$ws.Range("N13:O13").Merge()

# Fill in the new row's data - "CONGESTAL 20 TABS"
Set-TextValue 13 3 "CONGESTAL 20 TABS"
Set-TextValue 13 8 "0:1"
Set-TextValue 13 12 "1"
Set-TextValue 13 14 "50.00"
Set-TextValue 13 16 "25.0000"
Set-TextValue 13 17 "0:1"

# ------------------------------------------------------------------
# 2) Column A holds a plain running index (1..N); after the insert it
#    was shifted along with everything else, so re-number rows 7..47
#    sequentially.
# ------------------------------------------------------------------
$n = 1
for ($r = 7; $r -le 47; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $n
    $n++
}

# ------------------------------------------------------------------
# 3) Update the two cells that changed values (besides the shift):
#    VOLTAREN 75MG/3ML 3 AMP. row (now row 32) and the
#    "سرنجات 3 سم" row (now row 39).
# ------------------------------------------------------------------
Set-TextValue 32 8 "2:1"
Set-TextValue 32 16 "33.6600"
Set-TextValue 32 17 "0:2"

Set-TextValue 39 16 "6.0000"
Set-TextValue 39 17 "3:0"

# ------------------------------------------------------------------
# 4) Update the grand-total cell (now row 48) and the timestamp
#    footer (now row 49).
# ------------------------------------------------------------------
$ws.Cells.Item(48, 16).Value2 = 1581.8800000000001

Set-TextValue 49 1 "Saturday, 9 August, 2025 3:18 PM"
